# Updates cryptocurrency price/volume figures on Sheet1 (D/E columns).
# Values are textual (e.g. "307.08", "0.04%"), matching the source data
# feed which stores them as plain text, not numbers/percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2" = "307.08"
    "E2" = "0.04%"
    "D3" = "39.15"
    "E3" = "7.73%"
    "E4" = "0.55%"
    "D5" = "0.08049"
    "E5" = "-0.27%"
    "E6" = "2.34%"
    "D7" = "4.196"
    "E7" = "1.21%"
    "D8" = "7.960"
    "E8" = "2.18%"
    "D9" = "0.9319"
    "E9" = "0.69%"
    "D10" = "0.1456"
    "E10" = "1.02%"
    "D11" = "0.1926"
    "E11" = "0.76%"
    "D12" = "0.09056"
    "E12" = "0.10%"
    "D13" = "0.03510"
    "E13" = "2.18%"
    "D14" = "0.09789"
    "E14" = "-1.17%"
    "D15" = "0.001397"
    "E15" = "-1.04%"
    "D16" = "0.005884"
    "E16" = "-2.62%"
    "D17" = "3.794"
    "E17" = "-1.18%"
    "D18" = "3.415"
    "E18" = "0.68%"
    "D19" = "0.3446"
    "E19" = "-0.16%"
    "D20" = "0.1303"
    "E20" = "-2.42%"
    "D21" = "4.783"
    "E21" = "-0.36%"
    "D22" = "0.2508"
    "E22" = "-4.12%"
    "D23" = "0.04371"
    "E23" = "-0.13%"
    "D24" = "0.001238"
    "E24" = "0.59%"
    "D25" = "0.004276"
    "E25" = "-0.39%"
    "D39" = "0.02048"
    "E39" = "1.92%"
    "D40" = "0.05042"
    "E40" = "-1.95%"
    "D41" = "0.007449"
    "E41" = "-0.46%"
    "D42" = "0.01012"
    "E42" = "-0.04%"
    "D43" = "0.1352"
    "E43" = "-0.49%"
    "D44" = "0.002140"
    "E44" = "-0.32%"
    "D45" = "0.009056"
    "E45" = "-6.37%"
    "E46" = "-1.62%"
    "E47" = "0.32%"
    "D49" = "0.001602"
    "E49" = "28.09%"
    "D50" = "0.00002103"
    "E50" = "0.32%"
    "D51" = "0.0002003"
    "E51" = "0.32%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force the value to be stored as text (not auto-converted to a
    # number/percentage) by writing it while number-formatted as Text,
    # then clearing the format again so no style change is left behind.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.ClearFormats()
}
